# Apply cryptos list update (prices / 1h volume %) per commit
# "Updated cryptos list on Thu Mar  7 04:27:45 UTC 2024 with GitHub Actions"
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "65.959.11"
$ws.Range("E2").Value = "  +3.74%  "
# Row 3
$ws.Range("D3").Value = "3.793.23"
$ws.Range("E3").Value = "  +6.57%  "
# Row 4
$ws.Range("E4").Value = "  -0.23%  "
# Row 5
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "427.76"
$ws.Range("D5").ClearFormats()
$ws.Range("E5").Value = "  +7.80%  "
# Row 6
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "138.63"
$ws.Range("D6").ClearFormats()
$ws.Range("E6").Value = "  +12.74%  "
# Row 7
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.621"
$ws.Range("D7").ClearFormats()
$ws.Range("E7").Value = "  +5.71%  "
# Row 8
$ws.Range("E8").Value = "  -0.11%  "
# Row 9
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.740"
$ws.Range("D9").ClearFormats()
$ws.Range("E9").Value = "  +8.70%  "
# Row 10
$ws.Range("E10").Value = "  +0.75%  "
# Row 11
$ws.Range("E11").Value = "  -3.93%  "
# Row 12
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "42.89"
$ws.Range("D12").ClearFormats()
$ws.Range("E12").Value = "  +10.28%  "
# Row 13
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "10.53"
$ws.Range("D13").ClearFormats()
$ws.Range("E13").Value = "  +14.76%  "
# Row 14
$ws.Range("D14").Value = "4.394.22"
$ws.Range("E14").Value = "  +6.68%  "
# Row 15
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "15.09"
$ws.Range("D15").ClearFormats()
$ws.Range("E15").Value = "  +9.42%  "
# Row 16
$ws.Range("B16").Value = "WrappedEther"
$ws.Range("C16").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D16").Value = "3.797.36"
$ws.Range("E16").Value = "  +7.22%  "
# Row 17
$ws.Range("B17").Value = "TRON"
$ws.Range("C17").Value = "https://coinranking.com/coin/qUhEFk1I61atv+tron-trx"
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "0.138"
$ws.Range("D17").ClearFormats()
$ws.Range("E17").Value = "  +1.19%  "
# Row 18
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "20.01"
$ws.Range("D18").ClearFormats()
$ws.Range("E18").Value = "  +7.01%  "
# Row 19
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "1.13"
$ws.Range("D19").ClearFormats()
$ws.Range("E19").Value = "  +10.93%  "
# Row 20
$ws.Range("D20").Value = "66.167.19"
$ws.Range("E20").Value = "  +3.96%  "
# Row 21
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "406.00"
$ws.Range("D21").ClearFormats()
$ws.Range("E21").Value = "  +2.97%  "
# Row 22
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "15.17"
$ws.Range("D22").ClearFormats()
$ws.Range("E22").Value = "  +8.99%  "
# Row 23
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "3.24"
$ws.Range("D23").ClearFormats()
$ws.Range("E23").Value = "  +11.21%  "
# Row 24
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "85.21"
$ws.Range("D24").ClearFormats()
$ws.Range("E24").Value = "  +4.07%  "
# Row 25
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "36.77"
$ws.Range("D25").ClearFormats()
$ws.Range("E25").Value = "  +8.43%  "
# Row 26
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "9.96"
$ws.Range("D26").ClearFormats()
$ws.Range("E26").Value = "  +47.16%  "
# Row 27
$ws.Range("E27").Value = "  +9.53%  "
# Row 28
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "9.90"
$ws.Range("D28").ClearFormats()
$ws.Range("E28").Value = "  +13.22%  "
# Row 29
$ws.Range("E29").Value = "  -0.49%  "
# Row 30
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "13.91"
$ws.Range("D30").ClearFormats()
$ws.Range("E30").Value = "  +16.38%  "
# Row 31
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "704.40"
$ws.Range("D31").ClearFormats()
$ws.Range("E31").Value = "  +3.82%  "
# Row 32
$ws.Range("E32").Value = "  +16.71%  "
# Row 33
$ws.Range("E33").Value = "  +7.33%  "
# Row 34
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "40.72"
$ws.Range("D34").ClearFormats()
$ws.Range("E34").Value = "  +11.05%  "
# Row 35
$ws.Range("E35").Value = "  +0.05%  "
# Row 36
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "5.75"
$ws.Range("D36").ClearFormats()
$ws.Range("E36").Value = "  +41.45%  "
# Row 37
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "0.151"
$ws.Range("D37").ClearFormats()
# Row 38
$ws.Range("E38").Value = "  +5.27%  "
# Row 39
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "0.0476"
$ws.Range("D39").ClearFormats()
$ws.Range("E39").Value = "  +8.72%  "
# Row 40
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "2.67"
$ws.Range("D40").ClearFormats()
$ws.Range("E40").Value = "  +50.12%  "
# Row 41
$ws.Range("B41").Value = "PEPE"
$ws.Range("C41").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D41").Value = "0.0₃0682"
$ws.Range("E41").Value = "  +3.12%  "
# Row 42
$ws.Range("B42").Value = "Stellar"
$ws.Range("C42").Value = "https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm"
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "0.142"
$ws.Range("D42").ClearFormats()
$ws.Range("E42").Value = "  +8.78%  "
# Row 43
$ws.Range("B43").Value = "ThetaToken"
$ws.Range("C43").Value = "https://coinranking.com/coin/B42IRxNtoYmwK+thetatoken-theta"
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "2.86"
$ws.Range("D43").ClearFormats()
$ws.Range("E43").Value = "  +6.81%  "
# Row 44
$ws.Range("E44").Value = "  +0.31%  "
# Row 45
$ws.Range("E45").Value = "  +9.61%  "
# Row 46
$ws.Range("E46").Value = "  +16.34%  "
# Row 47
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "3.14"
$ws.Range("D47").ClearFormats()
$ws.Range("E47").Value = "  +1.89%  "
# Row 48
$ws.Range("E48").Value = "  +7.36%  "
# Row 49
$ws.Range("E49").Value = "  +5.73%  "
# Row 50
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "142.35"
$ws.Range("D50").ClearFormats()
$ws.Range("E50").Value = "  +0.44%  "
# Row 51
$ws.Range("E51").Value = "  +6.26%  "
